$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new baseline counter cell A3 = 0 ---
$ws.Range("A3").Value2 = 0

# --- Column B: copy the student name (currently only in column C) into column B ---
# so that column A can compare B (entered name) against C (roster name).
for ($r = 4; $r -le 33; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 3).Value2
}

# --- Column M:P (rows 4-33): rebuild as one shared formula group (matches si="1") ---
$ws.Range("M4:P33").Formula = '=IF($L4=M$2,1,0)'

# --- Column A row 4: first (non-shared) running-count formula ---
$ws.Range("A4").Formula = '=IF(B4=C4,1+A3,"_______")'

# --- Column A rows 5-34: shared formula group (matches si="2", ref A5:A34) ---
$ws.Range("A5:A34").Formula = '=IF(B5=C5,1+A4,"_______")'

# Row 34 itself keeps no formula/value, just inherits the same cell format as A4.
$ws.Range("A34").ClearContents()
$ws.Range("A4").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 35 grew slightly taller after the extra bottom border landed on row 34 ---
$ws.Rows.Item(35).RowHeight = 13

# --- View state: selection / scroll position ---
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 3
$ws.Range("B34").Select()
